$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-09-19 Friday" "2025-09-20 Saturday"

Replace-Text "170×9=" "936×7="
Replace-Text "827×7=" "184×8="
Replace-Text "387×7=" "864×3="
Replace-Text "106×3=" "808×9="
Replace-Text "753×4=" "681×5="
Replace-Text "326×9=" "382×2="
Replace-Text "292×4=" "181×2="
Replace-Text "779×4=" "337×6="
Replace-Text "215×5=" "401×7="
Replace-Text "425×9=" "549×2="
Replace-Text "994×5=" "657×7="
Replace-Text "950×4=" "760×7="
Replace-Text "424×2=" "446×9="
Replace-Text "472×2=" "933×6="
Replace-Text "168×3=" "357×7="
Replace-Text "676×7=" "931×2="
Replace-Text "932×6=" "385×6="
Replace-Text "103×5=" "766×5="
Replace-Text "206×9=" "477×7="
Replace-Text "353×6=" "287×4="
Replace-Text "261×9=" "943×4="
Replace-Text "702×2=" "974×7="
Replace-Text "983×4=" "687×9="
Replace-Text "333×9=" "446×3="
Replace-Text "244×4=" "830×6="
